# Weekly refresh of the "Poroto granado" price sheet:
# insert a brand-new record as row 10 (pushing the existing rows 10-48
# down to 11-49) so the data stays sorted with the newest observation on top.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 10:48 down to 11:49, leaving row 10 blank.
$ws.Rows("10:10").Insert()

# Populate the new row 10 with the latest market observation.
$ws.Range("A10").Value = 10
$ws.Range("B10").Value = "Vega Modelo de Temuco"
$ws.Range("C10").Value = "La Araucanía"
$ws.Range("D10").Value = 44558
$ws.Range("E10").Value = 9
$ws.Range("F10").Value = 100112030
$ws.Range("G10").Value = "Poroto granado"
$ws.Range("H10").Value = "Sin especificar"
$ws.Range("I10").Value = "Primera"
$ws.Range("J10").Value = 28
$ws.Range("K10").Value = 30000
$ws.Range("L10").Value = 30000
$ws.Range("M10").Value = 30000
$ws.Range("N10").Value = "$/saco 25 kilos"
$ws.Range("O10").Value = "Provincia de Limarí"
$ws.Range("P10").Value = 1200
$ws.Range("Q10").Value = 25
$ws.Range("R10").Value = "Hortaliza"
